# scritp que busca patrones por dia seleccionado
#
# Agrega a la hoja "Hoja1" los datos historicos de los dias nuevos
# (19 columnas nuevas, GD:GV) para las filas 2 a 7, y deja seleccionado
# el bloque recien incorporado (GT2:GV7), igual que hace Excel al
# terminar de pegar/escribir datos en la ultima columna usada.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Una fila por cada fila de datos (2..7), una columna por cada nueva
# columna agregada (GD=186 .. GV=204).
$data = @(
    @(4,0,0,0,12,1,6,0,1,0,20,20,3,14,0,0,7,0,7),
    @(10,4,2,7,13,9,7,10,9,8,22,25,18,20,0,11,8,3,17),
    @(19,9,11,15,15,16,8,11,14,24,24,33,10,26,22,14,16,8,23),
    @(25,13,16,13,23,28,20,13,18,25,25,34,20,27,23,23,26,11,28),
    @(33,26,17,31,28,29,21,32,30,33,27,35,22,29,33,29,29,24,31),
    @(34,34,22,34,20,23,22,35,33,34,36,36,30,30,36,32,35,36,35)
)

$firstRow = 2
$firstCol = 186   # GD

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $firstRow + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $firstCol + $j
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}

# Selecciona el rango con las columnas nuevas (coincide con lo que
# queda seleccionado tras cargar el dia agregado mas reciente).
$null = $ws.Range("GT2:GV7").Select()
